# "Generate Report for Handoff"
#
# The localization-status workbook tracks, per target language, whether a
# file is still "In Translation" or has been packaged up and is
# "Ready for handoff" — plus the timestamp of the most recent handoff
# generation run. This run just finished, so:
#   - flip every "In Translation" status cell to "Ready for handoff"
#   - bump the "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     timestamps to the new run time
#   - widen the Status column(s) so the longer "Ready for handoff" label
#     isn't clipped

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ------------------
$wsOverview.Range("E2").Value = "Ready for handoff"   # zh-cn status column
$wsOverview.Range("F2").Value = "Ready for handoff"   # de-de status column
$wsZhCn.Range("C2").Value     = "Ready for handoff"   # zh-cn table Status
$wsDeDe.Range("C2").Value     = "Ready for handoff"   # de-de table Status

# --- Latest handoff timestamps -----------------------------------------
$wsOverview.Range("G2").Value = "2016-08-20 12:43:08" # Latest HO Xliff Generate Date
$wsDeDe.Range("H2").Value     = "2016-08-20 12:43:08" # de-de Latest Handoff Datetime
$wsZhCn.Range("H2").Value     = "2016-08-20 12:43:01" # zh-cn Latest Handoff Datetime

# --- Widen the Status column(s) to fit "Ready for handoff" -------------
$wsOverview.Columns.Item(5).ColumnWidth = 16.33  # column E (zh-cn status)
$wsOverview.Columns.Item(6).ColumnWidth = 16.33  # column F (de-de status)
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.33  # column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.33  # column C (Status)
